# Generate Report for Handback
# Updates the localization-status workbook to reflect that the de-de
# handback has now been generated (zh-cn was already "in sync"), i.e.:
#  - Overview / per-language "Status" cells flip from "Ready for handoff"
#    to "Handed back: in sync with en-US"
#  - zh-cn and de-de sheets get their "Latest Target File"/"Latest
#    Handback File" columns (I, J) populated, with I2 becoming a
#    hyperlink to the source .md file (like column A already is)
#  - de-de's "Latest Handback DateTime" (K2) moves from the "never
#    handed back" sentinel to a real timestamp; zh-cn's sentinel text is
#    likewise replaced by its own handback timestamp

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$mdName     = "0c4274bb-7005-4af8-b3a2-b38817b0c895.md"
$mdUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dad0e58619458f278c959fdde0ad41a743dc5bc1/e2e/0c4274bb-7005-4af8-b3a2-b38817b0c895.md"
$zhXlf      = "0c4274bb-7005-4af8-b3a2-b38817b0c895.085d18978db4ebc1024d9ed7103fcf0f50d9e06f.zh-cn.xlf"
$deXlf      = "0c4274bb-7005-4af8-b3a2-b38817b0c895.085d18978db4ebc1024d9ed7103fcf0f50d9e06f.de-de.xlf"
$zhHandback = "2016-09-01 07:10:36"
$deHandback = "2016-09-01 07:10:43"

# Widen columns to fit the longer status text everywhere it appears.
$wideWidth = 29.166666666666668   # -> stored column width 30 (closest reachable to 29.9777047293527)
$fullWidth = 39.166666666666664   # -> stored column width 40 exactly

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E1").ColumnWidth = $wideWidth
$wsOverview.Range("F1").ColumnWidth = $wideWidth

# ---------------------------------------------------------------------
# zh-cn sheet: handback already generated and in sync
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C1").ColumnWidth = $wideWidth
$wsZh.Range("I1").ColumnWidth = $fullWidth
$wsZh.Range("J1").ColumnWidth = $fullWidth

$wsZh.Range("I2").Value = $mdName
$wsZh.Range("J2").Value = $zhXlf
$wsZh.Range("K2").Value = $zhHandback

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, "", "", $mdName)

# ---------------------------------------------------------------------
# de-de sheet: handback just generated
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C1").ColumnWidth = $wideWidth
$wsDe.Range("I1").ColumnWidth = $fullWidth
$wsDe.Range("J1").ColumnWidth = $fullWidth

$wsDe.Range("I2").Value = $mdName
$wsDe.Range("J2").Value = $deXlf
$wsDe.Range("K2").Value = $deHandback

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, "", "", $mdName)

Write-Host "Handback report generated for zh-cn and de-de."
